# "Corrected code of Parser" - fill in missing parser-table entries and
# append new states (rows 35-38) to the LR parsing table on Sheet1.
#
# Note: new shared-strings must be introduced in the same order they first
# appear in the target workbook's sharedStrings table (R2, S34, R14, S36),
# so the writes below are ordered accordingly rather than strictly by
# worksheet position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- introduce the four brand-new shared strings in target order ---------
$ws.Range("I35").Value = "R2"    # new string #1
$ws.Range("E9").Value  = "S34"   # new string #2
$ws.Range("F38").Value = "R14"   # new string #3
$ws.Range("F37").Value = "S36"   # new string #4

# --- remaining edits to existing rows -------------------------------------
$ws.Range("B18").Value = "S23"
$ws.Range("U18").Value = 22
$ws.Range("U19").Value = 25
$ws.Range("N27").Value = "R11"
$ws.Range("I33").Value = "R6"

# --- row 35 (new) ----------------------------------------------------------
$ws.Range("A35").Value = 33
$ws.Range("M35").Value = "R2"

# --- row 36 (new) ----------------------------------------------------------
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "S12"
$ws.Range("E36").Value = "S34"
$ws.Range("S36").Value = 35
$ws.Range("T36").Value = 10
$ws.Range("U36").Value = 11

# --- row 37 (new) ------------------------------------------------------
$ws.Range("A37").Value = 35
$ws.Range("H37").Value = "S16"
$ws.Range("J37").Value = "S15"

# --- row 38 (new) ------------------------------------------------------
$ws.Range("A38").Value = 36
$ws.Range("H38").Value = "R14"
$ws.Range("J38").Value = "R14"
$ws.Range("K38").Value = "R14"
$ws.Range("L38").Value = "R14"
$ws.Range("N38").Value = "R14"

# --- view state: move the active cell/selection the way the author left it -
$ws.Range("F37").Select() | Out-Null
